$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2056.8948
$ws.Range("I33").Value = 1498.3846
$ws.Range("J33").Value = 3267
$ws.Range("K33").Value = 1498.3846
$ws.Range("L33").Value = 3267
$ws.Range("M33").Value = -1269.3846
$ws.Range("N33").Value = -3725

$ws.Range("H62").Value = 51475956
$ws.Range("I62").Value = 21745054
$ws.Range("J62").Value = 113640570
$ws.Range("K62").Value = 21745054
$ws.Range("L62").Value = 113640570
$ws.Range("M62").Value = -21744430
$ws.Range("N62").Value = -113641818

$ws.Range("H65").Value = 51475956
$ws.Range("I65").Value = 21745054
$ws.Range("J65").Value = 113640570
$ws.Range("K65").Value = 108725270
$ws.Range("L65").Value = 568202850
$ws.Range("M65").Value = -108722150
$ws.Range("N65").Value = -568209090

$ws.Range("H98").Value = 37458790
$ws.Range("I98").Value = 28574944
$ws.Range("J98").Value = 42641036
$ws.Range("K98").Value = 28574944
$ws.Range("L98").Value = 42641036
$ws.Range("M98").Value = -28573446
$ws.Range("N98").Value = -42644032

$ws.Range("H113").Value = 2858804.2
$ws.Range("I113").Value = 4168189.5
$ws.Range("J113").Value = 1964
$ws.Range("K113").Value = 4168189.5
$ws.Range("L113").Value = 1964
$ws.Range("M113").Value = -4164935.5
$ws.Range("N113").Value = -8472

$ws.Range("H116").Value = 12969308
$ws.Range("J116").Value = 16676217
$ws.Range("L116").Value = 16676217
$ws.Range("N116").Value = -16683101

$ws.Range("H122").Value = 37458790
$ws.Range("I122").Value = 28574944
$ws.Range("J122").Value = 42641036
$ws.Range("K122").Value = 85724832
$ws.Range("L122").Value = 127923108
$ws.Range("M122").Value = -85722382
$ws.Range("N122").Value = -127928008

$ws.Range("H132").Value = 7938619.5
$ws.Range("I132").Value = 2456.4
$ws.Range("K132").Value = 7369.200000000001
$ws.Range("M132").Value = -4839.200000000001

$ws.Range("H137").Value = 19750552
$ws.Range("I137").Value = 4311151
$ws.Range("J137").Value = 69499736
$ws.Range("K137").Value = 12933453
$ws.Range("L137").Value = 208499208
$ws.Range("M137").Value = -12930903
$ws.Range("N137").Value = -208504308


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6654555.5
$ws.Range("I32").Value = 1607012
$ws.Range("J32").Value = 31261332
$ws.Range("K32").Value = 1607012
$ws.Range("L32").Value = 31261332
$ws.Range("M32").Value = -1606725
$ws.Range("N32").Value = -31261906

$ws.Range("H122").Value = 1293.7
$ws.Range("I122").Value = 603.3333
$ws.Range("J122").Value = 1858.5454
$ws.Range("K122").Value = 1809.9999
$ws.Range("L122").Value = 5575.6362
$ws.Range("M122").Value = 640.0001
$ws.Range("N122").Value = -10475.6362


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1960.1428
$ws.Range("I105").Value = 2043.75
$ws.Range("J105").Value = 1848.6666
$ws.Range("K105").Value = 2043.75
$ws.Range("L105").Value = 1848.6666
$ws.Range("M105").Value = -296.75
$ws.Range("N105").Value = -5342.6666


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2744864.2
$ws.Range("I31").Value = 1544823.1
$ws.Range("J31").Value = 5690419.5
$ws.Range("K31").Value = 1544823.1
$ws.Range("L31").Value = 5690419.5
$ws.Range("M31").Value = -1544528.1
$ws.Range("N31").Value = -5691009.5

$ws.Range("H34").Value = 2744864.2
$ws.Range("I34").Value = 1544823.1
$ws.Range("J34").Value = 5690419.5
$ws.Range("K34").Value = 1544823.1
$ws.Range("L34").Value = 5690419.5
$ws.Range("M34").Value = -1544621.1
$ws.Range("N34").Value = -5690823.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 940.11365
$ws.Range("J131").Value = 986.6
$ws.Range("L131").Value = 2959.8
$ws.Range("N131").Value = -13039.8

$ws.Range("H140").Value = 2777.6775
$ws.Range("I140").Value = 2706.4
$ws.Range("J140").Value = 2907.2727
$ws.Range("K140").Value = 8119.200000000001
$ws.Range("L140").Value = 8721.8181
$ws.Range("M140").Value = -2939.200000000001
$ws.Range("N140").Value = -19081.8181


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2669043.2
$ws.Range("I70").Value = 1592952.5
$ws.Range("K70").Value = 1592952.5
$ws.Range("M70").Value = -1592682.5

$ws.Range("H73").Value = 2669043.2
$ws.Range("I73").Value = 1592952.5
$ws.Range("K73").Value = 1592952.5
$ws.Range("M73").Value = -1592016.5

$ws.Range("H102").Value = 7428
$ws.Range("I102").Value = 8657.538
$ws.Range("K102").Value = 8657.538
$ws.Range("M102").Value = -7035.538

$ws.Range("H122").Value = 11942894
$ws.Range("I122").Value = 58150.11
$ws.Range("J122").Value = 33335434
$ws.Range("K122").Value = 174450.33
$ws.Range("L122").Value = 100006302
$ws.Range("M122").Value = -172000.33
$ws.Range("N122").Value = -100011202

$ws.Range("H126").Value = 8612
$ws.Range("I126").Value = 12283.777
$ws.Range("J126").Value = 2002.8
$ws.Range("K126").Value = 36851.331
$ws.Range("L126").Value = 6008.4
$ws.Range("M126").Value = -34381.331
$ws.Range("N126").Value = -10948.4

$ws.Range("H132").Value = 17819230
$ws.Range("I132").Value = 22511706
$ws.Range("J132").Value = 11367075
$ws.Range("K132").Value = 67535118
$ws.Range("L132").Value = 34101225
$ws.Range("M132").Value = -67532588
$ws.Range("N132").Value = -34106285


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2783.2856
$ws.Range("I40").Value = 999
$ws.Range("J40").Value = 3080.6667
$ws.Range("K40").Value = 999
$ws.Range("L40").Value = 3080.6667
$ws.Range("M40").Value = -863
$ws.Range("N40").Value = -3352.6667

$ws.Range("H122").Value = 9222365
$ws.Range("I122").Value = 1184542.1
$ws.Range("K122").Value = 3553626.3
$ws.Range("M122").Value = -3551176.3


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1230.3636
$ws.Range("I122").Value = 1150.9
$ws.Range("J122").Value = 1352.6154
$ws.Range("K122").Value = 3452.7
$ws.Range("L122").Value = 4057.8462
$ws.Range("M122").Value = -1002.7
$ws.Range("N122").Value = -8957.8462

$ws.Range("H126").Value = 41669828
$ws.Range("I126").Value = 83333860
$ws.Range("J126").Value = 5799.6665
$ws.Range("K126").Value = 250001580
$ws.Range("L126").Value = 17398.9995
$ws.Range("M126").Value = -249999110
$ws.Range("N126").Value = -22338.9995

